# Team_Everyday_Attendence.xlsx - add 3 new attendance rows (31-Aug / "No
# Meeting" separator / 04-Sep / 05-Sep) plus the reviewer comments that went
# with the 04-Sep & 05-Sep rows, matching commit:
#   "B8-B1=Team Attendance - 05-Sep-23"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 28: a lone "No Meeting" marker row (no date, no other columns)
# ---------------------------------------------------------------------
$ws.Range("A28").Value = "No Meeting"

# ---------------------------------------------------------------------
# Row 29: Thursday 31-Aug-2023
# ---------------------------------------------------------------------
$ws.Range("A29").NumberFormat = $ws.Range("A27").NumberFormat
$ws.Range("A29").Value = "8/31/2023"
$ws.Range("B29").Value = "PRESENT"
$ws.Range("C29").Value = "PRESENT"
$ws.Range("D29").Value = "PRESENT"
$ws.Range("E29").Value = "PRESENT"
$ws.Range("F29").Value = "ABSENT"
$ws.Range("G29").Value = "ABSENT"
$ws.Range("H29").Value = "ABSENT"
$ws.Range("I29").Value = "ABSENT"
$ws.Range("J29").Value = "ABSENT"
$ws.Range("K29").Value = "ABSENT"

# ---------------------------------------------------------------------
# Row 30: Monday 04-Sep-2023
# ---------------------------------------------------------------------
$ws.Range("A30").NumberFormat = $ws.Range("A27").NumberFormat
$ws.Range("A30").Value = "9/4/2023"
$ws.Range("B30").Value = "PRESENT"
$ws.Range("C30").Value = "PRESENT"
$ws.Range("D30").Value = "PRESENT"
$ws.Range("E30").Value = "PRESENT"
$ws.Range("F30").Value = "PRESENT"
$ws.Range("G30").Value = "ABSENT"
$ws.Range("H30").Value = "ABSENT"
$ws.Range("I30").Value = "ABSENT"
$ws.Range("J30").Value = "ABSENT"
$ws.Range("K30").Value = "ABSENT"

# ---------------------------------------------------------------------
# Row 31: Tuesday 05-Sep-2023
# ---------------------------------------------------------------------
$ws.Range("A31").NumberFormat = $ws.Range("A27").NumberFormat
$ws.Range("A31").Value = "9/5/2023"
$ws.Range("B31").Value = "PRESENT"
$ws.Range("C31").Value = "PRESENT"
$ws.Range("D31").Value = "PRESENT"
$ws.Range("E31").Value = "PRESENT"
$ws.Range("F31").Value = "PRESENT"
$ws.Range("G31").Value = "ABSENT"
$ws.Range("H31").Value = "PRESENT"
$ws.Range("I31").Value = "ABSENT"
$ws.Range("J31").Value = "PRESENT"
$ws.Range("K31").Value = "ABSENT"

# ---------------------------------------------------------------------
# Reviewer comments ("LENOVO") explaining the ABSENT marks above
# ---------------------------------------------------------------------
$ws.Range("G30").AddComment("LENOVO:" + [char]10 + "No Response")
$ws.Range("I30").AddComment("LENOVO:" + [char]10 + "No Response")
$ws.Range("J30").AddComment("LENOVO:" + [char]10 + "No Response")
$ws.Range("K30").AddComment("LENOVO:" + [char]10 + "No Response")
$ws.Range("G31").AddComment("LENOVO:" + [char]10 + "No response")
$ws.Range("K31").AddComment("LENOVO:" + [char]10 + "No response")

# ---------------------------------------------------------------------
# Leave the selection where the author left it
# ---------------------------------------------------------------------
$ws.Range("H31").Select()
